# "Added trees in Gui"
# Insert a new "SEPERATOR" label into column A of every blank separator row
# (the style-only rows, style index 3, that previously had no text) on the
# ReportNew sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReportNew")

$separatorRows = @(
    4, 11, 18, 25, 34, 69, 81, 88, 95, 104, 111, 120, 131, 137, 160, 165,
    171, 176, 180, 185, 191, 197, 203, 225, 229, 233, 237, 241, 245, 249,
    253, 257, 261, 265, 269, 273, 277
)

foreach ($r in $separatorRows) {
    $ws.Cells.Item($r, 1).Value = "SEPERATOR"
}

# Reflect the minimized window state recorded in the saved workbook view.
$excel.WindowState = -4140
